$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update forecast coal outage figures (MW)
$ws.Range("C5").Value = 840
$ws.Range("D5").Value = 1560
$ws.Range("E5").Value = 1560
$ws.Range("F5").Value = 1560
$ws.Range("G5").Value = 1560
$ws.Range("H5").Value = 1910

# Row 11: update expected return date and days until return for Gladstone 4
$ws.Range("H11").Value = 45738
$ws.Range("I11").Value = 8
